# Weekly update: insert a new price-report row for "Palta" (avocado) as row 42,
# pushing the existing rows 42..81 down to 43..82 (dimension grows to A1:T82).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("42:42").Insert()

$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C42").Value = "Arica y Parinacota"
$ws.Range("D42").Value = 44629
$ws.Range("E42").Value = 15
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100106
$ws.Range("H42").Value = "Oleaginosos"
$ws.Range("I42").Value = 100106002
$ws.Range("J42").Value = "Palta"
$ws.Range("K42").Value = "Hass"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 74000
$ws.Range("O42").Value = 75000
$ws.Range("P42").Value = 74500
$ws.Range("Q42").Value = "`$/caja 25 kilos"
$ws.Range("R42").Value = "Región de Coquimbo"
$ws.Range("S42").Value = 2980
$ws.Range("T42").Value = 25
